$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 227.75
$ws.Range("I28").Value = 227.75
$ws.Range("K28").Value = 227.75
$ws.Range("M28").Value = 257.25
$ws.Range("H74").Value = 2874.5
$ws.Range("I74").Value = 2711
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2711
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1775
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 2874.5
$ws.Range("I77").Value = 2711
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 13555
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -8875
$ws.Range("N77").Value = -34360
$ws.Range("H94").Value = 1075
$ws.Range("I94").Value = 1075
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1075
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -624
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 10875.8
$ws.Range("I113").Value = 6715.75
$ws.Range("K113").Value = 6715.75
$ws.Range("M113").Value = -3461.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 890.5
$ws.Range("I2").Value = 822.7778
$ws.Range("K2").Value = 822.7778
$ws.Range("M2").Value = -709.7778
$ws.Range("H74").Value = 820.5333000000001
$ws.Range("I74").Value = 639.46155
$ws.Range("J74").Value = 1997.5
$ws.Range("K74").Value = 639.46155
$ws.Range("L74").Value = 1997.5
$ws.Range("M74").Value = 234.53845
$ws.Range("N74").Value = -3745.5
$ws.Range("H77").Value = 820.5333000000001
$ws.Range("I77").Value = 639.46155
$ws.Range("J77").Value = 1997.5
$ws.Range("K77").Value = 3197.30775
$ws.Range("L77").Value = 9987.5
$ws.Range("M77").Value = 1170.69225
$ws.Range("N77").Value = -18723.5
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
$ws.Range("H96").Value = 19333
$ws.Range("J96").Value = 19333
$ws.Range("L96").Value = 19333
$ws.Range("N96").Value = -24825
$ws.Range("H116").Value = 890.5
$ws.Range("I116").Value = 822.7778
$ws.Range("K116").Value = 822.7778
$ws.Range("M116").Value = 1471.2222
$ws.Range("H132").Value = 1458.4286
$ws.Range("I132").Value = 1431.35
$ws.Range("K132").Value = 4294.049999999999
$ws.Range("M132").Value = -1764.049999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 890.5
$ws.Range("I3").Value = 822.7778
$ws.Range("K3").Value = 822.7778
$ws.Range("M3").Value = -708.7778
$ws.Range("H54").Value = 14297.5
$ws.Range("I54").Value = 14297.5
$ws.Range("K54").Value = 14297.5
$ws.Range("M54").Value = -13813.5
$ws.Range("H75").Value = 130000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 130000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 130000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -131872
$ws.Range("H78").Value = 130000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 130000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 390000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -399360
$ws.Range("H99").Value = 2715.5
$ws.Range("I99").Value = 2458.6
$ws.Range("K99").Value = 2458.6
$ws.Range("M99").Value = -960.5999999999999
$ws.Range("H105").Value = 3344.3
$ws.Range("I105").Value = 3160.4443
$ws.Range("J105").Value = 4999
$ws.Range("K105").Value = 3160.4443
$ws.Range("L105").Value = 4999
$ws.Range("M105").Value = -1413.4443
$ws.Range("N105").Value = -8493
$ws.Range("H130").Value = 110780
$ws.Range("J130").Value = 110780
$ws.Range("L130").Value = 110780
$ws.Range("N130").Value = -120820
$ws.Range("H134").Value = 11269.637
$ws.Range("I134").Value = 11329.667
$ws.Range("J134").Value = 10999.5
$ws.Range("K134").Value = 33989.001
$ws.Range("L134").Value = 32998.5
$ws.Range("M134").Value = -31454.001
$ws.Range("N134").Value = -38068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2572.75
$ws.Range("I31").Value = 1791.5385
$ws.Range("K31").Value = 1791.5385
$ws.Range("M31").Value = -1496.5385
$ws.Range("H34").Value = 2572.75
$ws.Range("I34").Value = 1791.5385
$ws.Range("K34").Value = 1791.5385
$ws.Range("M34").Value = -1589.5385
$ws.Range("H94").Value = 2298.3333
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H132").Value = 3366.3333
$ws.Range("I132").Value = 3049.5
$ws.Range("K132").Value = 9148.5
$ws.Range("M132").Value = -6618.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1482
$ws.Range("J75").Value = 502.5
$ws.Range("L75").Value = 1507.5
$ws.Range("N75").Value = -3503.5
$ws.Range("H78").Value = 1482
$ws.Range("J78").Value = 502.5
$ws.Range("L78").Value = 4522.5
$ws.Range("N78").Value = -14506.5
$ws.Range("H99").Value = 1400
$ws.Range("I99").Value = 1400
$ws.Range("K99").Value = 4200
$ws.Range("M99").Value = -1954
$ws.Range("H107").Value = 733.3333
$ws.Range("I107").Value = 600.75
$ws.Range("K107").Value = 1802.25
$ws.Range("M107").Value = 117.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 39500
$ws.Range("J47").Value = 39500
$ws.Range("L47").Value = 39500
$ws.Range("N47").Value = -40636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20320.723
$ws.Range("I7").Value = 20320.723
$ws.Range("K7").Value = 20320.723
$ws.Range("M7").Value = -20208.723
$ws.Range("H22").Value = 809.7778
$ws.Range("I22").Value = 572.4666999999999
$ws.Range("K22").Value = 572.4666999999999
$ws.Range("M22").Value = -277.4666999999999
$ws.Range("H27").Value = 809.7778
$ws.Range("I27").Value = 572.4666999999999
$ws.Range("K27").Value = 572.4666999999999
$ws.Range("M27").Value = -465.4666999999999
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H126").Value = 20320.723
$ws.Range("I126").Value = 20320.723
$ws.Range("K126").Value = 60962.16900000001
$ws.Range("M126").Value = -58492.16900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2274.5
$ws.Range("I96").Value = 2050.75
$ws.Range("J96").Value = 2722
$ws.Range("K96").Value = 2050.75
$ws.Range("L96").Value = 2722
$ws.Range("M96").Value = -677.75
$ws.Range("N96").Value = -5468
$ws.Range("H100").Value = 979.2
$ws.Range("I100").Value = 1080.5
$ws.Range("K100").Value = 2161
$ws.Range("M100").Value = -1620
$ws.Range("H101").Value = 21666.334
$ws.Range("J101").Value = 21666.334
$ws.Range("L101").Value = 21666.334
$ws.Range("N101").Value = -28156.334
$ws.Range("H132").Value = 3244.5557
$ws.Range("I132").Value = 3174.5
$ws.Range("K132").Value = 9523.5
$ws.Range("M132").Value = -6993.5
